# Fruta / hortaliza, semanal
# A new daily price-report row for "Mango" at Vega Central Mapocho de Santiago
# was inserted into the weekly consolidation, right before the existing row
# for date 44575 (previously row 307). Every subsequent row shifts down by
# one (old row 434 becomes row 435), and the sheet's used range grows by a
# row accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 307; everything from the old row 307 onward
# (through the old row 434) shifts down to 308..435.
$ws.Rows.Item(307).Insert()

# Populate the newly inserted row with the new day's record.
$ws.Range("A307").Value = 9
$ws.Range("B307").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C307").Value = "Metropolitana"
$ws.Range("D307").Value = 44726
$ws.Range("E307").Value = 13
$ws.Range("F307").Value = "Fruta"
$ws.Range("G307").Value = 100108
$ws.Range("H307").Value = "Tropicales y subtropicales"
$ws.Range("I307").Value = 100108002
$ws.Range("J307").Value = "Mango"
$ws.Range("K307").Value = "Sin especificar"
$ws.Range("L307").Value = "Primera"
$ws.Range("M307").Value = 630
$ws.Range("N307").Value = 7500
$ws.Range("O307").Value = 8500
$ws.Range("P307").Value = 7944
$ws.Range("Q307").Value = "$/bandeja 4 kilos"
$ws.Range("R307").Value = "Brasil"
$ws.Range("S307").Value = 1986
$ws.Range("T307").Value = 4
